$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("T3").Value = 0.99
$ws.Range("U3").Value = 0.99
$ws.Range("V3").Value = 0.99
$ws.Range("W3").Value = 0.99
